# Oferta.xlsx update: refresh date, replace product list with new items,
# remove product photo, extend totals table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the offer date header (C1)
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "Oferta z dnia: 15.10.2023"

# ---------------------------------------------------------------------
# 2. Remove the product photo (the logo "Picture 1" stays)
# ---------------------------------------------------------------------
$ws.Shapes.Item("Picture 2").Delete()

# ---------------------------------------------------------------------
# 3. Update the first product row (row 4) with the new first product
# ---------------------------------------------------------------------
$ws.Range("B4").Value = "A01655"
$ws.Range("C4").Value = "Gwiazdka wisząca NATURE_Aluro L " + [char]10 + "Szerokość: 23,00 cm " + [char]10 + "Głębokość: 4,00 cm " + [char]10 + "Wysokość: 23,00 cm " + [char]10
$ws.Range("D4").Value = "11,80 zł"
$ws.Range("E4").Value = "1/6 szt."
$ws.Range("G4").Formula = "=D4*F4"

# Row 4 no longer is the last product row, so its bottom border becomes thin
# (matching the other internal table rows) instead of the heavy closing rule.
$ws.Range("A4:F4").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------
# 4. Insert three more rows before the totals row for the extra products
# ---------------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Row 5 - second product (same look as row 4)
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$ws.Rows.Item(5).RowHeight = 130
$ws.Range("B5").Value = "A01656"
$ws.Range("C5").Value = "Gwiazdka wisząca NATURE_Aluro XL " + [char]10 + "Szerokość: 32,00 cm " + [char]10 + "Głębokość: 5,00 cm " + [char]10 + "Wysokość: 32,00 cm " + [char]10
$ws.Range("D5").Value = "17,90 zł"
$ws.Range("E5").Value = "1/4 szt."
$ws.Range("G5").Formula = "=D5*F5"

# Row 6 - third product (same look as row 4)
$ws.Range("A4:G4").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)
$ws.Rows.Item(6).RowHeight = 130
$ws.Range("B6").Value = "A01670"
$ws.Range("C6").Value = "Renifer stojący_Aluro XL " + [char]10 + "Szerokość: 32,00 cm " + [char]10 + "Głębokość: 5,00 cm " + [char]10 + "Wysokość: 29,00 cm " + [char]10
$ws.Range("D6").Value = "76,00 zł"
$ws.Range("E6").Value = "1/8 szt."
$ws.Range("G6").Formula = "=D6*F6"

# Row 7 - fourth (last) product, closes the table with a thick bottom rule
$ws.Range("A4:G4").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)
$ws.Rows.Item(7).RowHeight = 130
$ws.Range("B7").Value = "A01667"
$ws.Range("C7").Value = "Renifer świecznik_Aluro " + [char]10 + "Szerokość: 22,00 cm " + [char]10 + "Głębokość: 10,00 cm " + [char]10 + "Wysokość: 23,00 cm " + [char]10
$ws.Range("D7").Value = "32,90 zł"
$ws.Range("E7").Value = "1/4 szt."
$ws.Range("G7").Formula = "=D7*F7"
$ws.Range("A7:F7").Borders.Item(9).Weight = 4

# ---------------------------------------------------------------------
# 5. Fix up the totals row (pushed from row 5 down to row 8)
# ---------------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 30
$ws.Range("G8").Formula = "=SUM(G4:G7)"

Write-Host "Edit applied"
